# Update cfb_weather.xlsx with Timestamp 2024-09-22T10:01:16.328967

$wb = $excel.ActiveWorkbook

# --- Update the shared Timestamp string on the "FBS" sheet (column AD, every data row) ---
$wsFbs = $wb.Worksheets.Item("FBS")
$oldTimestamp = "2024-09-22T05:15:50.836917"
$newTimestamp = "2024-09-22T10:01:16.328967"

$usedRange = $wsFbs.UsedRange
$lastRow = $usedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsFbs.Cells.Item($r, 30)  # column AD
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}

# --- Update the weather figures on the "Other" sheet, row 2 ---
$wsOther = $wb.Worksheets.Item("Other")
$wsOther.Range("K2").Value = 79.45999999999999
$wsOther.Range("L2").Value = 7.2
$wsOther.Range("M2").Value = 0
$wsOther.Range("P2").Value = -6.8
